# Moved some configuration settings to config and added support for
# multiple sequence passes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink the header row now that some configuration text moved out of it.
$ws.Rows.Item(1).RowHeight = 81.75

# D3 now references the new "2,3" multi-pass sequence value (shared string)
# instead of the literal number 2.
$ws.Range("D3").Value = "2,3"

# B4 toggles from "N" to "Y" now that this test case participates too.
$ws.Range("B4").Value = "Y"

# New D4 cell carries the extra sequence-pass count for this row.
$ws.Range("D4").Value = 3

# Update the current selection to match the last-edited cell.
$ws.Range("B5").Select()
